$d = $word.ActiveDocument

# --- First paragraph (the **ID__...__ID** placeholder line) ---
$p = $d.Paragraphs(1)

# Add a paragraph border (top/left/bottom/right) with 5pt space-from-text,
# no explicit line (matches <w:pBdr><w:top w:space="5"/>...</w:pBdr>)
$b = $p.Borders
$b.DistanceFromTop = 5
$b.DistanceFromLeft = 5
$b.DistanceFromBottom = 5
$b.DistanceFromRight = 5

# Change left indent from 120 twips (6pt) to 225 twips (11.25pt)
$p.LeftIndent = 11.25

# Replace the placeholder id text and drop the trailing " " run that
# followed it (collapsing the paragraph down to a single run).
$oldId = "**ID__AFFARS_pgi_5343_topic_4__ID**"
$newId = "**ID__AFFARS_SMC_PGI_5343__ID**"

$idRange = $d.Range(0, $oldId.Length)
$idRange.Text = $newId

$spaceStart = $newId.Length
$spaceRange = $d.Range($spaceStart, $spaceStart + 1)
if ($spaceRange.Text -eq " ") {
    $spaceRange.Delete()
}
